$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Bob Johnson/Data Structures/C -> John Doe/Introduction to Databases/A
$ws.Range("A12").Value = "John Doe"
$ws.Range("B12").Value = "Introduction to Databases"
$ws.Range("C12").Value = "A"

# Row 13: Bob Johnson/Data Structures/C -> John Doe/Web Development/B
$ws.Range("A13").Value = "John Doe"
$ws.Range("B13").Value = "Web Development"
$ws.Range("C13").Value = "B"

# Rows 14-16 (Bob Johnson/Data Structures/C) remain unchanged

# Row 17: Jane Smith/Introduction to Databases/A -> Bob Johnson/Data Structures/C
$ws.Range("A17").Value = "Bob Johnson"
$ws.Range("B17").Value = "Data Structures"
$ws.Range("C17").Value = "C"

# Row 18: Jane Smith/Introduction to Databases/A -> Bob Johnson/Data Structures/C
$ws.Range("A18").Value = "Bob Johnson"
$ws.Range("B18").Value = "Data Structures"
$ws.Range("C18").Value = "C"

# Row 19: Jane Smith/Introduction to Databases/A -> Bob Johnson/Data Structures/C
$ws.Range("A19").Value = "Bob Johnson"
$ws.Range("B19").Value = "Data Structures"
$ws.Range("C19").Value = "C"

# Rows 20-21 (Jane Smith/Introduction to Databases/A) remain unchanged

# New rows 22-25: Jane Smith/Introduction to Databases/A
$ws.Range("A22").Value = "Jane Smith"
$ws.Range("B22").Value = "Introduction to Databases"
$ws.Range("C22").Value = "A"

$ws.Range("A23").Value = "Jane Smith"
$ws.Range("B23").Value = "Introduction to Databases"
$ws.Range("C23").Value = "A"

$ws.Range("A24").Value = "Jane Smith"
$ws.Range("B24").Value = "Introduction to Databases"
$ws.Range("C24").Value = "A"

$ws.Range("A25").Value = "Jane Smith"
$ws.Range("B25").Value = "Introduction to Databases"
$ws.Range("C25").Value = "A"
